$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "tool_checkout_log" — sign-out/sign-in rework
# ------------------------------------------------------------------
$log = $wb.Worksheets.Item("tool_checkout_log")

# Drop the old scratch/demo rows (rows 5-8) before rewriting the top rows,
# so the row-2..4 values below land on a clean sheet.
$log.Rows("5:8").Delete()

# Headers
$log.Range("A1").Value = "Sign Out Employee"
$log.Range("B1").Value = "Tool"
$log.Range("C1").Value = "Sign Out Time"
$log.Range("D1").Value = "Sign In Time"
$log.Range("E1").Value = "Sign In Employee"

# New checkout/checkin data
$log.Range("A2").Value = "emp5"
$log.Range("B2").Value = "tool5"
$log.Range("C2").Value = "02/10/2024 14:03"
$log.Range("D2").Value = "02/10/2024 14:05"
$log.Range("E2").Value = "emp1"

$log.Range("A3").Value = "emp3"
$log.Range("B3").Value = "tool3"
$log.Range("C3").Value = "02/10/2024 14:04"
$log.Range("D3").Value = "02/10/2024 14:05"
$log.Range("E3").Value = "emp2"

$log.Range("A4").Value = "emp6"
$log.Range("B4").Value = "tool6"
$log.Range("C4").Value = "02/10/2024 14:04"
$log.Range("D4").Value = "02/10/2024 14:05"
$log.Range("E4").Value = "emp6"

# Column widths (characters). Excel's ColumnWidth COM setter snaps to the
# nearest whole pixel (1/6-character) of the workbook's Normal-style font,
# so the inputs below are chosen to land on the closest achievable stored
# width to the target (15.33203125 / 15.5 / 17.6640625 / 19.1640625).
$log.Columns.Item(1).ColumnWidth = 14.5
$log.Columns.Item(3).ColumnWidth = 14.666667
$log.Columns.Item(4).ColumnWidth = 16.833333
$log.Columns.Item(5).ColumnWidth = 18.333333

$log.Range("F16").Select() | Out-Null

# ------------------------------------------------------------------
# Sheet "employees" — trim roster back to 10 entries
# ------------------------------------------------------------------
$employees = $wb.Worksheets.Item("employees")
$employees.Rows("12:16").Delete()
$employees.Columns.Item(1).ColumnWidth = 14.666667
$employees.Range("F13").Select() | Out-Null

# ------------------------------------------------------------------
# Sheet "tools" — trim catalog back to 10 entries
# ------------------------------------------------------------------
$tools = $wb.Worksheets.Item("tools")
$tools.Rows("12:16").Delete()
$tools.Range("D8").Select() | Out-Null

$log.Activate()
